$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-40, replacing the old Strike#-derived
# counts with the regenerated K counts.
$gValues = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 5
    12 = 5
    13 = 2
    14 = 1
    15 = 1
    16 = 2
    17 = 2
    18 = 3
    19 = 1
    20 = 2
    21 = 0
    22 = 3
    23 = 0
    24 = 1
    25 = 1
    26 = 2
    27 = 0
    28 = 1
    29 = 3
    30 = 2
    31 = 2
    32 = 1
    33 = 2
    34 = 1
    35 = 2
    36 = 3
    37 = 1
    38 = 3
    39 = 1
    40 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
